$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "sigmazoid"
$ws.Range("A5").Value = "Rokas Miceika"
$ws.Range("B5").Value = "miceikarokas24@gmail.com"

$table = $ws.ListObjects.Item("user")
$table.Resize($ws.Range("A1:B5"))
